# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Adds a new period (2509) row to the "Estado de Cuenta" table, updates the
# totals (Valor Mora / Cant. Periodos) accordingly, and keeps the signature
# block (rows that used to be 25-26) as the last two rows of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row for period 2509 right after the current last
#        data row (20), pushing the blank spacer rows and the signature
#        block (old rows 21-26) down by one. -----------------------------
$ws.Rows("21").Insert()

# The newly inserted row 21 is currently blank; row 20 still carries the
# "last row of the table" border/format. Re-home the formats: row 21
# (the new last row) gets row 20's old "closing" borders, and row 20
# reverts to the regular "middle of table" borders (copied from row 19).
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Center the "Periodo Mora" column for the whole data block.
$ws.Range("E16:E21").HorizontalAlignment = -4108   # xlCenter

# --- 2. Fill in the data for the new period 2509, matching the existing
#        worker row pattern. ----------------------------------------------
$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1034916526"
$ws.Range("D21").Value2 = "MARIA ALEJANDRA HERNANDEZ ORTIZ"
$ws.Range("E21").Value2 = "2509"
$ws.Range("F21").Value2 = 60960
$ws.Range("G21").Value2 = 1524000

# --- 3. Update the summary figures: one more overdue period is now
#        included, so "Valor Mora" and "Cant. Periodos" both grow. --------
$ws.Range("E11").Value2 = 339344
$ws.Range("F13").Value2 = 6
